$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The last data row (row 18, the "float" switch-combination row) needs to move up
# to become the new row 6, with the existing rows 6-17 shifting down by one
# row each to become rows 7-18.
#
# Implement as: insert a blank row above row 6 (pushing old rows 6-18 down to
# 7-19), copy the row that is now at 19 (originally row 18) into the new
# row 6, then delete the now-duplicated row 19.
$ws.Rows.Item(6).Insert(-4121) | Out-Null
$ws.Range("A19:G19").Copy($ws.Range("A6:G6")) | Out-Null
$ws.Rows.Item(19).Delete() | Out-Null

# Update the active cell selection shown in the sheet view.
$ws.Activate()
$ws.Range("I12").Select() | Out-Null
